$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while we write the new price strings, so that
# values such as "59.408.36" or "540.45" are stored as text, matching the source data.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "59.408.36"
$ws.Range("E2").Value = "  -4.28%  "
$ws.Range("D3").Value = "2.481.37"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "540.45"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").Value = "146.47"
$ws.Range("E6").Value = "  -4.97%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").Value = "2.500.18"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "2.907.88"
$ws.Range("E14").Value = "  -5.22%  "
$ws.Range("D15").Value = "24.27"
$ws.Range("E15").Value = "  -5.01%  "
$ws.Range("D16").Value = "59.228.27"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "2.495.60"
$ws.Range("E18").Value = "  -4.35%  "
$ws.Range("D19").Value = "11.48"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "4.37"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").Value = "326.11"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").Value = "0.994"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "5.79"
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("D24").Value = "61.12"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  -10.44%  "
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "0.991"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "2.594.69"
$ws.Range("E28").Value = "  -4.92%  "
$ws.Range("D29").Value = "7.90"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").Value = "7.15"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "0.0₃0788"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").Value = "0.995"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "159.23"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").Value = "1.43"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").Value = "18.69"
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("D38").Value = "4.48"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").Value = "5.97"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "317.64"
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("D42").Value = "36.55"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("D43").Value = "3.76"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").Value = "0.839"
$ws.Range("E44").Value = "  -5.53%  "
$ws.Range("D45").Value = "0.993"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "0.595"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").Value = "10.74"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").Value = "125.57"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "0.0936"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").Value = "0.0527"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "0.0231"
$ws.Range("E51").Value = "  -3.40%  "

# Restore the original (default) cell style on column D so formatting matches the source file.
$dRange.Style = "Normal"

Write-Host "Updated cryptos list"
